$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: replace numeric zip codes with quoted-text zip codes -------
# (values were re-entered from a source that carried them as quoted strings,
#  e.g. "97212" instead of the bare number 97212)
$ws.Range("D2").Value = '"97212"'
$ws.Range("D3").Value = '"97213"'
$ws.Range("D4").Value = '"97214"'
$ws.Range("D5").Value = '"97209"'
$ws.Range("D6").Value = '"97203"'
$ws.Range("D7").Value = '"97217"'
$ws.Range("D8").Value = '"97210"'
$ws.Range("D9").Value = '"97212"'
$ws.Range("D10").Value = '"97211"'
$ws.Range("D11").Value = '"97232"'
$ws.Range("D12").Value = '"97214"'
$ws.Range("D13").Value = '"97202"'
$ws.Range("D14").Value = '"97206"'
$ws.Range("D15").Value = '"97202"'
$ws.Range("D16").Value = '"97217"'
$ws.Range("D17").Value = '"97266"'
$ws.Range("D18").Value = '"97060"'
$ws.Range("D19").Value = '"97203"'
$ws.Range("D20").Value = '"97213"'
$ws.Range("D21").Value = '"97206"'
$ws.Range("D22").Value = '"97266"'
$ws.Range("D23").Value = '"97220"'
$ws.Range("D24").Value = '"97030"'
$ws.Range("D25").Value = '"97217"'
$ws.Range("D26").Value = '"97201"'
$ws.Range("D27").Value = '"97232"'
$ws.Range("D28").Value = '"97203"'
$ws.Range("D29").Value = '"97217"'
$ws.Range("D30").Value = '"97209"'
$ws.Range("D31").Value = '"97232"'
$ws.Range("D32").Value = '"97219"'
$ws.Range("D33").Value = '"97214"'
$ws.Range("D34").Value = '"97213"'
$ws.Range("D35").Value = '"97220"'
$ws.Range("D36").Value = '"97236"'
$ws.Range("D37").Value = '"97060"'
$ws.Range("D38").Value = '"97203"'
$ws.Range("D39").Value = '"97211"'
$ws.Range("D40").Value = '"97213"'
$ws.Range("D41").Value = '"97232"'
$ws.Range("D42").Value = '"97209"'
$ws.Range("D43").Value = '"97201"'
$ws.Range("D44").Value = '"97214"'
$ws.Range("D45").Value = '"97202"'
$ws.Range("D46").Value = '"97206"'
$ws.Range("D47").Value = '"97219"'
$ws.Range("D48").Value = '"97236"'
$ws.Range("D49").Value = '"97080"'
$ws.Range("D50").Value = '"97030"'
$ws.Range("D51").Value = '"97060"'
$ws.Range("D52").Value = '"97217"'
$ws.Range("D53").Value = '"97205"'
$ws.Range("D54").Value = '"97220"'
$ws.Range("D55").Value = '"97202"'
$ws.Range("D56").Value = '"97216"'
$ws.Range("D57").Value = '"97024"'

# --- Column widths ----------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 68.83072916666667   # C -> 69.6640625 (xml chars)
$ws.Columns.Item(4).ColumnWidth = 14.666666666666666  # D -> 15.5 (xml chars)

# --- Number format on the income / population columns -----------------
# Drop the thousands-separator format (#,##0) in favor of a plain integer
# format (0) for the whole E:G columns (header + data).
$ws.Range("E1:G57").NumberFormat = "0"

# --- View / selection state --------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 75
$ws.Range("D3:D4").Select()
